$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column (H) with header text in H1 and a value in H2,
# matching the formatting used by the existing "sum" header in G1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H2").Value = 0
